$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the two naive-component forecast cells that were erroneously populated
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Correct tiny floating point recalculation values for the remaining forecasts
$ws.Range("C4").Value = 6.277541464866943
$ws.Range("C5").Value = 6.535114773304795
$ws.Range("E6").Value = 4.950888348161864
$ws.Range("E7").Value = 3.982564147794343
$ws.Range("E8").Value = 4.334309403335457
$ws.Range("E9").Value = 2.644356903452594
$ws.Range("E10").Value = 3.383932287548674
$ws.Range("C11").Value = 2.508469427909921
$ws.Range("E11").Value = 3.355044026998977
$ws.Range("E12").Value = 3.749841708124202
$ws.Range("C14").Value = 3.047037961814514
$ws.Range("E14").Value = 2.880436144359466
$ws.Range("C15").Value = -0.22288476972816
$ws.Range("C16").Value = -1.165854108406639
$ws.Range("E16").Value = 1.819118980963319
$ws.Range("E17").Value = 2.549024517027942
$ws.Range("C19").Value = 2.039329803030099
